# Refresh sync metadata (_uid / _updated / _updated-by) on several rows
# across sheets, drop the now-unused "_format" column from "Point Defs",
# and update a handful of data values.
# (commit: "can update all implemented types")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Defs"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Defs")

$ws.Range("A2").Value = "lgs4l68f-wupv"
$ws.Range("B2").Value = "2023-04-22T15:18:42.639Z"
$ws.Range("C2").Value = "lgs4l68f"

$ws.Range("A3").Value = "lgs4l68g-jlfh"
$ws.Range("B3").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C3").Value = "lgs4l68g"

$ws.Range("A4").Value = "lgs4l68g-095d"
$ws.Range("B4").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C4").Value = "lgs4l68g"

# ---------------------------------------------------------------------
# Sheet "Point Defs" - drop the "_format" column (L) entirely, refresh
# sync metadata, and update a couple of data values.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Point Defs")

$ws.Range("A2").Value = "lgs4l68g-0qsw"
$ws.Range("B2").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C2").Value = "lgs4l68g"
$ws.Range("F2").Value = "ats6"

$ws.Range("A3").Value = "lgs4l68g-mj7s"
$ws.Range("B3").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C3").Value = "lgs4l68g"

$ws.Range("A4").Value = "lgs4l68g-hvoj"
$ws.Range("B4").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C4").Value = "lgs4l68g"
$ws.Range("I4").Value = "Orig desc"

# Remove the whole "_format" column (L), shrinking the sheet to A1:K4
$ws.Range("L1:L4").EntireColumn.Delete()

# ---------------------------------------------------------------------
# Sheet "Entry Base"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Base")

$ws.Range("A2").Value = "lgs4l68g-0f7a"
$ws.Range("B2").Value = "2023-04-22T15:18:42.640Z"
$ws.Range("C2").Value = "lgs4l68g"
$ws.Range("G2").Value = "2023-04-22T06"
$ws.Range("H2").Value = "Orig note"

$ws.Range("A3").Value = "lgs4l68h-13pq"
$ws.Range("B3").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C3").Value = "lgs4l68h"
$ws.Range("F3").Value = "lgs4l68s-gttg"
$ws.Range("G3").Value = "2023-04-22T10:18:42"

# ---------------------------------------------------------------------
# Sheet "Entry Points"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Points")

$ws.Range("A2").Value = "lgs4l68h-w50n"
$ws.Range("B2").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C2").Value = "lgs4l68h"

$ws.Range("A3").Value = "lgs4l68h-1bns"
$ws.Range("B3").Value = "2023-04-22T15:18:42.641Z"
$ws.Range("C3").Value = "lgs4l68h"
